$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.772.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5329"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3733"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07140"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8861"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08130"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.902.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +37.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.285"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008486"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.69%  "

$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.812.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.969"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.380"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.278"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.739"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.692"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.620"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09097"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8069"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05022"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.169"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.949"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6091"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.683"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.178"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01939"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.065"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5275"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.478"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.751"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1486"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.645"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.962"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06061"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.05%  "

